$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.225.64'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.586.51'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.597.51'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('E11').Value = '  +3.10%  '
$ws.Range('E12').Value = '  +9.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.346'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.041.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.275.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.57'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.27%  '
$ws.Range('E17').Value = '  +4.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.590.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '337.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.68%  '
$ws.Range('E25').Value = '  +5.92%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('E28').Value = '  +2.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0782'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.15%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  +0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.07'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('E36').Value = '  +2.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.880'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.885'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.23%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '295.03'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.12%  '
$ws.Range('E42').Value = '  +1.77%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('E44').Value = '  +2.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.597'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.24%  '
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.08%  '
$ws.Range('E48').Value = '  -0.13%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.17%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0232'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.946.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.24%  '
